$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 986, shifting existing rows 986-1022 down to 989-1025
$ws.Rows.Item(986).Insert()
$ws.Rows.Item(987).Insert()
$ws.Rows.Item(988).Insert()

$ws.Cells.Item(986,1).Value2 = 7
$ws.Cells.Item(986,2).Value2 = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(986,3).Value2 = 'Ñuble'
$ws.Cells.Item(986,4).Value2 = 45075
$ws.Cells.Item(986,5).Value2 = 16
$ws.Cells.Item(986,6).Value2 = 100112004
$ws.Cells.Item(986,7).Value2 = 'Cebolla'
$ws.Cells.Item(986,8).Value2 = 'Morada(o)'
$ws.Cells.Item(986,9).Value2 = '1a (guarda)'
$ws.Cells.Item(986,10).Value2 = 80
$ws.Cells.Item(986,11).Value2 = 8000
$ws.Cells.Item(986,12).Value2 = 9000
$ws.Cells.Item(986,13).Value2 = 8375
$ws.Cells.Item(986,14).Value2 = '$/malla 18 kilos'
$ws.Cells.Item(986,15).Value2 = 'Región de O''Higgins'
$ws.Cells.Item(986,16).Value2 = 465
$ws.Cells.Item(986,17).Value2 = 18
$ws.Cells.Item(986,18).Value2 = 'Hortaliza'

$ws.Cells.Item(987,1).Value2 = 7
$ws.Cells.Item(987,2).Value2 = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(987,3).Value2 = 'Ñuble'
$ws.Cells.Item(987,4).Value2 = 45075
$ws.Cells.Item(987,5).Value2 = 16
$ws.Cells.Item(987,6).Value2 = 100112004
$ws.Cells.Item(987,7).Value2 = 'Cebolla'
$ws.Cells.Item(987,8).Value2 = 'Sin especificar'
$ws.Cells.Item(987,9).Value2 = '1a (guarda)'
$ws.Cells.Item(987,10).Value2 = 130
$ws.Cells.Item(987,11).Value2 = 7000
$ws.Cells.Item(987,12).Value2 = 8000
$ws.Cells.Item(987,13).Value2 = 7385
$ws.Cells.Item(987,14).Value2 = '$/malla 18 kilos'
$ws.Cells.Item(987,15).Value2 = 'Región de O''Higgins'
$ws.Cells.Item(987,16).Value2 = 410
$ws.Cells.Item(987,17).Value2 = 18
$ws.Cells.Item(987,18).Value2 = 'Hortaliza'

$ws.Cells.Item(988,1).Value2 = 7
$ws.Cells.Item(988,2).Value2 = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(988,3).Value2 = 'Ñuble'
$ws.Cells.Item(988,4).Value2 = 45075
$ws.Cells.Item(988,5).Value2 = 16
$ws.Cells.Item(988,6).Value2 = 100112004
$ws.Cells.Item(988,7).Value2 = 'Cebolla'
$ws.Cells.Item(988,8).Value2 = 'Sin especificar'
$ws.Cells.Item(988,9).Value2 = '1a (guarda)'
$ws.Cells.Item(988,10).Value2 = 100
$ws.Cells.Item(988,11).Value2 = 12000
$ws.Cells.Item(988,12).Value2 = 12000
$ws.Cells.Item(988,13).Value2 = 12000
$ws.Cells.Item(988,14).Value2 = '$/malla 25 kilos'
$ws.Cells.Item(988,15).Value2 = 'Región de O''Higgins'
$ws.Cells.Item(988,16).Value2 = 480
$ws.Cells.Item(988,17).Value2 = 25
$ws.Cells.Item(988,18).Value2 = 'Hortaliza'
